$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, $innerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $d.Paragraphs($paragraphIndex).Range.InsertXML($pkg)
}

$p1Xml = '<w:p w14:paraId="2A2A1502" w14:textId="3FF747B0" w:rsidR="008E31FB" w:rsidRDefault="00D938F0" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Diego Quan </w:t></w:r></w:p>'
Set-ParagraphXml 1 $p1Xml

$p2Xml = '<w:p w14:paraId="0C860142" w14:textId="2C4CCD44" w:rsidR="00D938F0" w:rsidRDefault="00D938F0" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>Victor Farfán</w:t></w:r></w:p>'
Set-ParagraphXml 2 $p2Xml

$p3Xml = '<w:p w14:paraId="196B7802" w14:textId="28D58596" w:rsidR="00D938F0" w:rsidRDefault="00D938F0" w:rsidP="00D938F0"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>COMPILADOR</w:t></w:r></w:p>'
Set-ParagraphXml 3 $p3Xml

$p5Xml = '<w:p w14:paraId="7E120F69" w14:textId="77777777" w:rsidR="00C7280A" w:rsidRDefault="00D938F0" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r w:rsidRPr="00D938F0"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>El archivo principal es el a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>rchivo “Compilador.py” en la carpeta compilers. Correr este archivo en la terminal por medio del comando “Compilador.py</w:t></w:r><w:r w:rsidR="00F36020"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>&lt;programa.decaf&gt;</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>” o “python3 Compilador.py</w:t></w:r><w:r w:rsidR="00F36020"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>&lt;programa.decaf&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="00F36020"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>&lt;programa.decaf&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> es un argumento obligatorio el cual contiene el programa input que va a ser compilado</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'
Set-ParagraphXml 5 $p5Xml

$p7Xml = '<w:p w14:paraId="173AC10D" w14:textId="77777777" w:rsidR="00C7280A" w:rsidRDefault="00D938F0" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Puede usar la bandera -target y asignarle el valor de “scan” para especificar la etapa del compilador a la que va a llegar el programa. Si no especifica la bandera -target, la etapa “scan” será elegida por defecto. </w:t></w:r></w:p>'
Set-ParagraphXml 7 $p7Xml

$p10Xml = '<w:p w14:paraId="083D38DA" w14:textId="581E183F" w:rsidR="00557085" w:rsidRDefault="00557085" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Para especificar un carácter debe usar comillas simples ‘ ’. </w:t></w:r></w:p>'
Set-ParagraphXml 10 $p10Xml

$p11Xml = '<w:p w14:paraId="6DED64ED" w14:textId="15951D2E" w:rsidR="00557085" w:rsidRPr="00D938F0" w:rsidRDefault="00557085" w:rsidP="00D938F0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-GT"/></w:rPr><w:t>Para especificar un String debe usar comillas dobles “ ”.</w:t></w:r></w:p>'
Set-ParagraphXml 11 $p11Xml

